$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values (E, G, H, I, J, K, M, N, O, P, Q, R, S, T) for rows 2-10
$data = @{
    2 = @{ E=3; G=7.235148; H=21.705444; I=0.9254344869740032; J=0.9254344869740032; K=3;
           M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118;
           Q=332.502177942356; R=2992.519601481204; S=0.5727688814317935; T=0.5727688814317935 }
    3 = @{ E=3; G=7.235148; H=21.705444; I=0.9254344869740032; J=0.9254344869740032; K=3;
           M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261;
           Q=49.560146400704; R=446.041317606336; S=0.08537240204919308; T=0.08537240204919307 }
    4 = @{ E=3; G=7.235148; H=21.705444; I=0.9254344869740032; J=0.9254344869740032; K=3;
           M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556;
           Q=155.168297705784; R=1396.514679352056; S=0.2672932034930167; T=0.2672932034930167 }
    5 = @{ E=3; G=0.3016356666666667; H=0.9049070000000001; I=0.03858166390441884; J=0.03858166390441884; K=3;
           M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118;
           Q=13.86212363752078; R=124.759112737687; S=0.02387892043073618; T=0.02387892043073618 }
    6 = @{ E=3; G=0.3016356666666667; H=0.9049070000000001; I=0.03858166390441884; J=0.03858166390441884; K=3;
           M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261;
           Q=2.066178577089778; R=18.595607193808; S=0.003559203129921193; T=0.003559203129921193 }
    7 = @{ E=3; G=0.3016356666666667; H=0.9049070000000001; I=0.03858166390441884; J=0.03858166390441884; K=3;
           M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556;
           Q=6.469016656468666; R=58.221149908218; S=0.01114354034376147; T=0.01114354034376147 }
    8 = @{ E=3; G=0.2813256666666666; H=0.843977; I=0.0359838491215779; J=0.0359838491215779; K=3;
           M=45.95651366666667; N=137.869541; O=0.6189188856627118; P=0.6189188856627118;
           Q=12.92874684495078; R=116.358721604557; S=0.02227108380018214; T=0.02227108380018214 }
    9 = @{ E=3; G=0.2813256666666666; H=0.843977; I=0.0359838491215779; J=0.0359838491215779; K=3;
           M=6.849914666666667; N=20.549744; O=0.09225115688993263; P=0.09225115688993261;
           Q=1.927056810209778; R=17.343511291888; S=0.003319551710818347; T=0.003319551710818346 }
    10 = @{ E=3; G=0.2813256666666666; H=0.843977; I=0.0359838491215779; J=0.0359838491215779; K=3;
            M=21.446458; N=64.33937399999999; O=0.2888299574473556; P=0.2888299574473556;
            Q=6.033439094488665; R=54.30095185039799; S=0.01039321361057741; T=0.01039321361057741 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
